$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.1824936751613834,
    0.2317473888535702,
    0.2376530296263142,
    0.1516086032871022,
    0.1990854234405351,
    0.2196813845232968,
    0.203489026958722,
    0.1485013260084578,
    0.1764744801046101,
    0.1794035352775858,
    0.1974074920200406,
    0.2102464003209671,
    0.1455353771332317,
    0.1799508157365124,
    0.1859361010674234,
    0.1829758758944881,
    0.2210789546456689,
    0.191290772714984,
    0.1676656629257355,
    0.1751824305185567,
    0.1738220095283137,
    0.2046199182240863,
    0.1158877922385684,
    0.1828970065198267,
    0.2669697286264353,
    0.1954848763931794,
    0.2136788875518469,
    0.1939492007172717,
    0.1847599881337664,
    0.1824931771449953,
    0.1885368251228011,
    0.1848221300218739,
    0.1583881705481297,
    0.1762584085916622,
    0.1733076713224236,
    0.1548617484913297,
    0.1584038075523559,
    0.2378545878153589,
    0.2475716845020496,
    0.2123844111806007,
    0.2277276277178376,
    0.1805944022943746,
    0.1480678615449562,
    0.2101712129795292,
    0.1793127956508198,
    0.2145078450610379,
    0.1860342742614958,
    0.2402290535652423,
    0.2365063749388666,
    0.2128173532484794,
    0.2304778916287789,
    0.1790222434576274,
    0.2096403935261369,
    0.2199922817222123,
    0.2204349446019456,
    0.1652770648031216,
    0.1863367111051463,
    0.1291021015699512,
    0.2049447649760349,
    0.2007773301522786,
    0.2384123800495567,
    0.2036775608132293,
    0.2476979212442918,
    0.2864219666562899,
    0.2442462490540736,
    0.2091784137032771,
    0.1845337331253921,
    0.1927607344470286,
    0.1696526789992195,
    0.1614397292013895,
    0.2130471213811032,
    0.1846785523602854,
    0.1883170164149853,
    0.2499334586778545,
    0.204048694661715,
    0.1922779627043707,
    0.149924146989616
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
